$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.997.88'
$ws.Range("D3").Value = '2.227.95'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.627'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.02'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.35%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0896'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.08%  '
$ws.Range("D13").Value = '2.564.79'
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.801'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").Value = '2.249.96'
$ws.Range("E18").Value = '  +2.46%  '
$ws.Range("D19").Value = '41.883.96'
$ws.Range("E19").Value = '  +4.54%  '
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0902'
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '168.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.73'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("E34").Value = '  +5.76%  '
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0635'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.64'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.42%  '
$ws.Range("E39").Value = '  -4.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000256'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +30.61%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("E42").Value = '  +4.57%  '
$ws.Range("E43").Value = '  -3.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.00%  '
$ws.Range("E45").Value = '  -0.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0968'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.71%  '
$ws.Range("D48").Value = '1.477.98'
$ws.Range("E48").Value = '  -2.77%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.91%  '
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.39%  '
